$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '71.552.47'
$ws.Range("E2").Value = '  -1.25%  '

# Row 3
$ws.Range("D3").Value = '3.879.71'
$ws.Range("E3").Value = '  -2.30%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.26'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +2.21%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.72'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  +9.11%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.672'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  -0.34%  '

# Row 8
$ws.Range("E8").Value = '  +0.30%  '

# Row 9
$ws.Range("E9").Value = '  +2.43%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.177'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  +6.53%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.08'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  +2.53%  '

# Row 12
$ws.Range("E12").Value = '  +1.79%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.43'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = '  +6.80%  '

# Row 14
$ws.Range("D14").Value = '4.505.64'
$ws.Range("E14").Value = '  -2.03%  '

# Row 15
$ws.Range("D15").Value = '3.887.96'
$ws.Range("E15").Value = '  -2.24%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.98'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  +3.13%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.93'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  +0.21%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.22'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  -4.04%  '

# Row 19
$ws.Range("E19").Value = '  -1.87%  '

# Row 20
$ws.Range("D20").Value = '71.431.96'
$ws.Range("E20").Value = '  -1.27%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '435.53'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +2.32%  '

# Row 22
$ws.Range("E22").Value = '  +1.81%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '94.29'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  -0.92%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.30'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  -3.53%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.86'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  -2.01%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.19'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  -4.33%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.00'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  -1.57%  '

# Row 28
$ws.Range("E28").Value = '  +0.12%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.21'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  -4.16%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.14'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  -2.69%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.05'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  +4.96%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '51.73'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  +5.34%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '13.73'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  +2.77%  '

# Row 34
$ws.Range("E34").Value = '  -3.08%  '

# Row 35
$ws.Range("D35").Value = '0.0₃0992'
$ws.Range("E35").Value = '  +16.92%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '68.44'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  +1.65%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '621.96'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  -8.07%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.424'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  -2.38%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  -0.03%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.31'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  +0.80%  '

# Row 41
$ws.Range("E41").Value = '  -0.06%  '

# Row 42
$ws.Range("E42").Value = '  -1.68%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.26'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  +43.65%  '

# Row 44
$ws.Range("E44").Value = '  -2.57%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.22'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  -5.94%  '

# Row 46
$ws.Range("B46").Value = 'Fetch.AI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.66'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  -1.97%  '

# Row 47
$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.145'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  -1.73%  '

# Row 48
$ws.Range("E48").Value = '  -0.61%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.81'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  -15.94%  '

# Row 50
$ws.Range("D50").Value = '2.855.98'
$ws.Range("E50").Value = '  +3.26%  '

# Row 51
$ws.Range("E51").Value = '  +1.41%  '
